$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ---- Sheet ALC ----
$ws1.Range("H5").Value = 0
$ws1.Range("J5").Value = 0
$ws1.Range("L5").Value = 0
$ws1.Range("N5").ClearContents()
$ws1.Range("H40").Value = 4241.4443
$ws1.Range("J40").Value = 3921.5
$ws1.Range("L40").Value = 3921.5
$ws1.Range("N40").Value = -4271.5
$ws1.Range("H98").Value = 3236.4
$ws1.Range("I98").Value = 3360.4285
$ws1.Range("K98").Value = 3360.4285
$ws1.Range("M98").Value = -1862.4285
$ws1.Range("H122").Value = 3236.4
$ws1.Range("I122").Value = 3360.4285
$ws1.Range("K122").Value = 10081.2855
$ws1.Range("M122").Value = -7631.2855
$ws1.Range("H138").Value = 3225.4363
$ws1.Range("J138").Value = 2953.8333
$ws1.Range("L138").Value = 8861.499899999999
$ws1.Range("N138").Value = -19141.4999

# ---- Sheet ARM ----
$ws2.Range("H2").Value = 568.4286
$ws2.Range("I2").Value = 291.09525
$ws2.Range("K2").Value = 291.09525
$ws2.Range("M2").Value = -178.09525
$ws2.Range("H32").Value = 5291.8945
$ws2.Range("I32").Value = 5406.6924
$ws2.Range("K32").Value = 5406.6924
$ws2.Range("M32").Value = -5119.6924
$ws2.Range("H61").Value = 2341.0278
$ws2.Range("I61").Value = 1540.48
$ws2.Range("J61").Value = 4160.4546
$ws2.Range("K61").Value = 1540.48
$ws2.Range("L61").Value = 4160.4546
$ws2.Range("M61").Value = -1328.48
$ws2.Range("N61").Value = -4584.4546
$ws2.Range("H103").Value = 40000
$ws2.Range("J103").Value = 40000
$ws2.Range("L103").Value = 40000
$ws2.Range("N103").Value = -42344
$ws2.Range("H116").Value = 568.4286
$ws2.Range("I116").Value = 291.09525
$ws2.Range("K116").Value = 291.09525
$ws2.Range("M116").Value = 2002.90475
$ws2.Range("H136").Value = 2341.0278
$ws2.Range("I136").Value = 1540.48
$ws2.Range("J136").Value = 4160.4546
$ws2.Range("K136").Value = 4621.440000000001
$ws2.Range("L136").Value = 12481.3638
$ws2.Range("M136").Value = -2071.440000000001
$ws2.Range("N136").Value = -17581.3638
$ws2.Range("H140").Value = 77842
$ws2.Range("J140").Value = 85815.664
$ws2.Range("L140").Value = 85815.664
$ws2.Range("N140").Value = -96175.664

# ---- Sheet BSM ----
$ws3.Range("H3").Value = 568.4286
$ws3.Range("I3").Value = 291.09525
$ws3.Range("K3").Value = 291.09525
$ws3.Range("M3").Value = -177.09525
$ws3.Range("H135").Value = 105239.5
$ws3.Range("J135").Value = 105239.5
$ws3.Range("L135").Value = 105239.5
$ws3.Range("N135").Value = -115379.5
$ws3.Range("H137").Value = 64999
$ws3.Range("J137").Value = 64999
$ws3.Range("L137").Value = 64999
$ws3.Range("N137").Value = -75199
$ws3.Range("H140").Value = 59533.637
$ws3.Range("J140").Value = 59533.637
$ws3.Range("L140").Value = 59533.637
$ws3.Range("N140").Value = -69893.637

# ---- Sheet CRP ----
$ws4.Range("H31").Value = 4965.457
$ws4.Range("I31").Value = 4054.2942
$ws4.Range("J31").Value = 5826
$ws4.Range("K31").Value = 4054.2942
$ws4.Range("L31").Value = 5826
$ws4.Range("M31").Value = -3759.2942
$ws4.Range("N31").Value = -6416
$ws4.Range("H34").Value = 4965.457
$ws4.Range("I34").Value = 4054.2942
$ws4.Range("J34").Value = 5826
$ws4.Range("K34").Value = 4054.2942
$ws4.Range("L34").Value = 5826
$ws4.Range("M34").Value = -3852.2942
$ws4.Range("N34").Value = -6230
$ws4.Range("J62").Value = 16872.75
$ws4.Range("L62").Value = 16872.75
$ws4.Range("N62").Value = -18120.75
$ws4.Range("J65").Value = 16872.75
$ws4.Range("L65").Value = 84363.75
$ws4.Range("N65").Value = -90603.75

# ---- Sheet CUL ----
$ws5.Range("H81").Value = 3967.3076
$ws5.Range("J81").Value = 4447.7
$ws5.Range("L81").Value = 13343.1
$ws5.Range("N81").Value = -15589.1
$ws5.Range("H84").Value = 3967.3076
$ws5.Range("J84").Value = 4447.7
$ws5.Range("L84").Value = 40029.3
$ws5.Range("N84").Value = -51261.3
$ws5.Range("H122").Value = 794.7778
$ws5.Range("J122").Value = 725.9167
$ws5.Range("L122").Value = 6533.2503
$ws5.Range("N122").Value = -11433.2503
$ws5.Range("H138").Value = 1580
$ws5.Range("I138").Value = 1580
$ws5.Range("J138").Value = 0
$ws5.Range("K138").Value = 4740
$ws5.Range("L138").Value = 0
$ws5.Range("M138").Value = 400
$ws5.Range("N138").ClearContents()
$ws5.Range("H141").Value = 18945.777
$ws5.Range("I141").Value = 8418.666999999999
$ws5.Range("K141").Value = 25256.001
$ws5.Range("M141").Value = -20076.001

# ---- Sheet GSM ----
$ws6.Range("H132").Value = 2815.4666
$ws6.Range("I132").Value = 2405.476
$ws6.Range("J132").Value = 3772.111
$ws6.Range("K132").Value = 7216.428
$ws6.Range("L132").Value = 11316.333
$ws6.Range("M132").Value = -4686.428
$ws6.Range("N132").Value = -16376.333
$ws6.Range("H138").Value = 98355.664
$ws6.Range("J138").Value = 98355.664
$ws6.Range("L138").Value = 98355.664
$ws6.Range("N138").Value = -108635.664

# ---- Sheet LTW ----
$ws7.Range("H68").Value = 3296.6
$ws7.Range("I68").Value = 3118
$ws7.Range("J68").Value = 3415.6667
$ws7.Range("K68").Value = 3118
$ws7.Range("L68").Value = 3415.6667
$ws7.Range("M68").Value = -2369
$ws7.Range("N68").Value = -4913.6667
$ws7.Range("H71").Value = 3296.6
$ws7.Range("I71").Value = 3118
$ws7.Range("J71").Value = 3415.6667
$ws7.Range("K71").Value = 15590
$ws7.Range("L71").Value = 17078.3335
$ws7.Range("M71").Value = -11846
$ws7.Range("N71").Value = -24566.3335
$ws7.Range("H74").Value = 49999
$ws7.Range("J74").Value = 49999
$ws7.Range("L74").Value = 49999
$ws7.Range("N74").Value = -51995
$ws7.Range("H77").Value = 49999
$ws7.Range("J77").Value = 49999
$ws7.Range("L77").Value = 149997
$ws7.Range("N77").Value = -159981
$ws7.Range("H82").Value = 2999.5
$ws7.Range("I82").Value = 2999.5
$ws7.Range("J82").Value = 0
$ws7.Range("K82").Value = 2999.5
$ws7.Range("L82").Value = 0
$ws7.Range("M82").Value = -2638.5
$ws7.Range("N82").ClearContents()
$ws7.Range("H85").Value = 2999.5
$ws7.Range("I85").Value = 2999.5
$ws7.Range("J85").Value = 0
$ws7.Range("K85").Value = 2999.5
$ws7.Range("L85").Value = 0
$ws7.Range("M85").Value = -1751.5
$ws7.Range("N85").ClearContents()

# ---- Sheet WVR ----
$ws8.Range("H136").Value = 83338150
$ws8.Range("I136").Value = 125000664
$ws8.Range("J136").Value = 13125
$ws8.Range("K136").Value = 375001992
$ws8.Range("L136").Value = 39375
$ws8.Range("M136").Value = -374999442
$ws8.Range("N136").Value = -44475
